# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 12:14:43"
$wsZhCn.Range("E3").Value = "2016-03-11 12:14:43"
$wsZhCn.Range("H2").Value = "2016-03-11 12:15:00"
$wsZhCn.Range("H3").Value = "2016-03-11 12:15:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 12:14:46"
$wsDeDe.Range("E3").Value = "2016-03-11 12:14:46"
$wsDeDe.Range("H2").Value = "2016-03-11 12:15:10"
$wsDeDe.Range("H3").Value = "2016-03-11 12:15:10"
